$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "乙肝表面抗原"
$ws.Range("A3").Value = "乙型肝炎表面抗体定量"
$ws.Range("A4").Value = "乙型肝炎e抗原定量"
$ws.Range("A5").Value = "乙型肝炎e抗体定量"

# B5 must remain text "2.7" (not a number) - use a leading apostrophe
# so Excel stores it as a text value with a quote prefix instead of
# coercing it into a numeric cell.
$ws.Range("B5").Value = "'2.7"

$ws.Range("A6").Value = "乙型肝炎核心抗体定量"
$ws.Range("A7").Value = "丙肝抗体"
$ws.Range("A8").Value = "丙肝核心抗原"
